$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 onto the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for rows 2-11
$saveValues = @(1, 1, 0, 1, 0, 0, 1, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
